# Add a new weekly price record for "Macroferia Regional de Talca - Zanahoria".
# A new row is inserted at row 553 (pushing the existing rows 553:578 down to
# 554:579). The new row duplicates the data that used to be in row 553,
# except for the date (column D) and the volume (column J), which get the
# new week's values: Fecha = 2023-08-09 (serial 45147), Volumen = 700.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$insertRow = 553
$lastCol = 18  # column R

# Capture the existing row 553 values before shifting anything down.
$existingValues = @()
for ($c = 1; $c -le $lastCol; $c++) {
    $existingValues += ,$ws.Cells.Item($insertRow, $c).Value2
}

# Shift rows 553:578 down to 554:579, preserving formatting.
$ws.Rows.Item($insertRow).Insert()

# Re-populate the freshly inserted row with a copy of the old row 553 data.
for ($c = 1; $c -le $lastCol; $c++) {
    $ws.Cells.Item($insertRow, $c).Value2 = $existingValues[$c - 1]
}

# Apply the new week's Fecha (D) and Volumen (J) values.
$ws.Cells.Item($insertRow, 4).Value2 = 45147
$ws.Cells.Item($insertRow, 10).Value2 = 700
